$wb = $excel.ActiveWorkbook

# --- Update selection on "CFRNostroInputter" (sheet 1) ---
$wsNostroInputter = $wb.Worksheets.Item(1)
$wsNostroInputter.Activate()
$wsNostroInputter.Range("H14").Select()

# --- Update selection on "CFRVostroInputter" (sheet 3) ---
$wsVostroInputter = $wb.Worksheets.Item(3)
$wsVostroInputter.Activate()
$wsVostroInputter.Range("L18").Select()

# --- Add the new "IBGCFRNostroInputter" sheet as a copy of "CFRNostroInputter", ---
# --- placed after the last existing sheet (so it becomes sheet 4) ---
$wsNostroInputter.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsIBG = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsIBG.Name = "IBGCFRNostroInputter"

# Make the new sheet the active tab/selection, matching the committed state
$wsIBG.Activate()
$wsIBG.Range("J17").Select()
